# Re-applies the latest scrape of cryptos.xlsx: refreshed Price (D) /
# Volume(1h) (E) figures for every coin row, a two-row reorder of
# ImmutableX <-> Filecoin (rows 30-31, incl. their Coin/Link/Price/Volume
# cells), and the knock-on Price/Volume tweaks around them.
#
# Column D holds plain text in the source sheet (e.g. "249.36"), not
# numbers -- many of the new quotes still look numeric, so a bare
# .Value assignment would let Excel's type-inference silently convert
# them to floats. Prefixing with a leading apostrophe forces literal
# text storage, exactly like typing '249.36 into the cell by hand,
# without touching the cell's (General) number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''37.092.41'
$ws.Range("E2").Value = '  +0.14%  '

# Row 3
$ws.Range("D3").Value = '''2.051.28'
$ws.Range("E3").Value = '  -0.57%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''249.36'
$ws.Range("E5").Value = '  -0.06%  '

# Row 6
$ws.Range("D6").Value = '''0.668'
$ws.Range("E6").Value = '  -0.92%  '

# Row 7
$ws.Range("D7").Value = '''59.35'
$ws.Range("E7").Value = '  +8.49%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("E9").Value = '  +1.36%  '

# Row 10
$ws.Range("D10").Value = '''0.0792'
$ws.Range("E10").Value = '  -0.36%  '

# Row 11
$ws.Range("E11").Value = '  +2.04%  '

# Row 12
$ws.Range("D12").Value = '''16.05'
$ws.Range("E12").Value = '  +6.89%  '

# Row 13
$ws.Range("D13").Value = '''2.350.49'
$ws.Range("E13").Value = '  -0.53%  '

# Row 14
$ws.Range("E14").Value = '  +2.40%  '

# Row 15
$ws.Range("D15").Value = '''5.77'
$ws.Range("E15").Value = '  +7.83%  '

# Row 16
$ws.Range("D16").Value = '''2.050.19'
$ws.Range("E16").Value = '  -0.67%  '

# Row 17
$ws.Range("D17").Value = '''18.29'
$ws.Range("E17").Value = '  +28.63%  '

# Row 18
$ws.Range("D18").Value = '''37.116.65'
$ws.Range("E18").Value = '  +0.22%  '

# Row 19
$ws.Range("D19").Value = '''75.76'
$ws.Range("E19").Value = '  +2.96%  '

# Row 20
$ws.Range("D20").Value = '''0.0₃0903'
$ws.Range("E20").Value = '  -3.29%  '

# Row 21
$ws.Range("E21").Value = '  +0.00%  '

# Row 22
$ws.Range("D22").Value = '''237.99'
$ws.Range("E22").Value = '  +0.21%  '

# Row 23
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("E24").Value = '  -0.39%  '

# Row 25
$ws.Range("D25").Value = '''2.21'
$ws.Range("E25").Value = '  +10.13%  '

# Row 26
$ws.Range("D26").Value = '''9.46'
$ws.Range("E26").Value = '  +4.28%  '

# Row 27
$ws.Range("D27").Value = '''169.07'
$ws.Range("E27").Value = '  -0.72%  '

# Row 28
$ws.Range("D28").Value = '''20.10'
$ws.Range("E28").Value = '  -0.09%  '

# Row 29
$ws.Range("E29").Value = '  +0.23%  '

# Row 30
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '''4.82'
$ws.Range("E30").Value = '  +4.30%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''1.12'
$ws.Range("E31").Value = '  +6.48%  '

# Row 32
$ws.Range("D32").Value = '''0.0627'
$ws.Range("E32").Value = '  -0.57%  '

# Row 33
$ws.Range("D33").Value = '''4.53'
$ws.Range("E33").Value = '  +2.79%  '

# Row 34
$ws.Range("D34").Value = '''0.0889'
$ws.Range("E34").Value = '  -0.11%  '

# Row 35
$ws.Range("E35").Value = '  -0.02%  '

# Row 36
$ws.Range("D36").Value = '''2.24'
$ws.Range("E36").Value = '  -2.43%  '

# Row 37
$ws.Range("E37").Value = '  -1.00%  '

# Row 38
$ws.Range("E38").Value = '  +4.30%  '

# Row 39
$ws.Range("E39").Value = '  -1.21%  '

# Row 40
$ws.Range("D40").Value = '''3.11'
$ws.Range("E40").Value = '  +11.13%  '

# Row 41
$ws.Range("D41").Value = '''5.15'
$ws.Range("E41").Value = '  +24.09%  '

# Row 42
$ws.Range("D42").Value = '''17.65'
$ws.Range("E42").Value = '  -1.21%  '

# Row 43
$ws.Range("D43").Value = '''0.0224'
$ws.Range("E43").Value = '  -0.82%  '

# Row 44
$ws.Range("E44").Value = '  -0.79%  '

# Row 45
$ws.Range("D45").Value = '''97.12'
$ws.Range("E45").Value = '  +0.06%  '

# Row 46
$ws.Range("E46").Value = '  +3.98%  '

# Row 47
$ws.Range("E47").Value = '  -5.18%  '

# Row 48
$ws.Range("D48").Value = '''1.289.34'
$ws.Range("E48").Value = '  -0.71%  '

# Row 49
$ws.Range("E49").Value = '  -1.32%  '

# Row 50
$ws.Range("D50").Value = '''6.82'
$ws.Range("E50").Value = '  -0.91%  '

# Row 51
$ws.Range("D51").Value = '''2.240.41'
$ws.Range("E51").Value = '  -0.35%  '
